# Applies the "Branch wise aging stock status" data realignment to the
# NoStock sheet.  The Item Name (D) / UOM (E) pairs for several rows were
# shuffled among each other (their shared-string table entries moved
# around). We reproduce the exact end state by writing each affected
# cell's final text directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dinafex group (rows 3-5) ---
$ws.Range("D3").Value = "Dinafex 120mg Tablet"
$ws.Range("D4").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 180mg Tablet"

# --- Etorix group (rows 7-9) ---
$ws.Range("D7").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E7").Value = "40's"
$ws.Range("D8").Value = "Etorix 90mg Tablet"
$ws.Range("E8").Value = "30's"
$ws.Range("D9").Value = "Etorix 120mg Tablet"
$ws.Range("E9").Value = "20's"

# --- Flucloxin group (rows 11-12) ---
$ws.Range("D11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E11").Value = "36 's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule"
$ws.Range("E12").Value = "30 's"

# --- Ketonic group (rows 14-16) ---
$ws.Range("D14").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E14").Value = "4's"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

# --- Kynol group (rows 19-20) UOM swap only ---
$ws.Range("E19").Value = "30 's"
$ws.Range("E20").Value = "30 's"

# --- Zithrox group (rows 25-27) ---
$ws.Range("D25").Value = "Zithrox 15ml Suspension"
$ws.Range("E25").Value = "15 ml"
$ws.Range("D26").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E26").Value = "30ml"
$ws.Range("D27").Value = "Zithrox 500mg Tablet"
$ws.Range("E27").Value = "6 's"
